$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.335.90'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '1.837.38'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("E4").Value = '  +1.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.94'
$ws.Range("E5").Value = '  +2.10%  '
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4738'
$ws.Range("E7").Value = '  +1.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3698'
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07460'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8844'
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.48'
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").Value = '1.876.83'
$ws.Range("E12").Value = '  +5.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07391'
$ws.Range("E13").Value = '  +4.66%  '
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.16'
$ws.Range("E15").Value = '  +1.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.576'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008837'
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.014'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.84'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").Value = '27.365.43'
$ws.Range("E21").Value = '  +1.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.354'
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.70'
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").Value = '2.078.29'
$ws.Range("E24").Value = '  +2.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.914'
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.26'
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.63'
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.170'
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.264'
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("E30").Value = '  +2.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08956'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7594'
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.177'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.556'
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.945'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.014'
$ws.Range("E36").Value = '  +1.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.105'
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05372'
$ws.Range("E38").Value = '  +1.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01961'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.003'
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.284'
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5350'
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1664'
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.534'
$ws.Range("E45").Value = '  +1.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4972'
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.50'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.015'
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.679'
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06320'
$ws.Range("E51").Value = '  +0.60%  '
